$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four oldest years (2006年-2009年); this shifts 2010年.. up by 4 rows.
$ws.Rows("2:5").Delete()

# After the shift, 2019年 is now on row 11 - update its changed figures.
$ws.Range("C11").Value = 7073.9
$ws.Range("D11").Value = 30.5

# Append the new 2021年 row (row 13), matching the formatting used for the
# other year-label cells in column A (bold, centered, thin border) by
# copying the format from the cell directly above.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"
$ws.Range("C13").Value = 6718.86
$ws.Range("D13").Value = 33.11
$ws.Range("E13").Value = 19634.49
$ws.Range("F13").Value = 220.63
$ws.Range("G13").Value = 20423.18
$ws.Range("H13").Value = 44.63
